# Generate Report for Handback
#
# Populates the "Latest Target File", "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns for row 8
# (the a14ebcf6-3fca-4c76-925c-380c384446cc entry) on both the
# zh-cn and de-de sheets, and widens the "Error Detail" column so the
# new message is readable.

$wb = $excel.ActiveWorkbook

$targetMdDisplay = "a14ebcf6-3fca-4c76-925c-380c384446cc.md"
$targetMdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f10f2c6f2549ec613174dd7dc62bcd3149096a19/e2e/a14ebcf6-3fca-4c76-925c-380c384446cc.md"
$errorDetail     = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7266f0a7e49197aad8a41a52fecf17e7cc2af314/e2e/a14ebcf6-3fca-4c76-925c-380c384446cc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f10f2c6f2549ec613174dd7dc62bcd3149096a19/e2e/a14ebcf6-3fca-4c76-925c-380c384446cc.md."

function Update-HandbackRow($sheetName, $handbackFile, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen the "Error Detail" column (P) so the long message fits.
    $ws.Range("P1").ColumnWidth = 39.17

    # I8: Latest Target File -> the handed-back markdown file, shown as a hyperlink
    # just like the other hyperlinked cells in this sheet.
    $ws.Range("I8").Value = $targetMdDisplay
    $ws.Range("I8").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null

    # J8: Latest Handback File
    $ws.Range("J8").Value = $handbackFile

    # K8: Latest Handback DateTime
    $ws.Range("K8").Value = $handbackDateTime

    # P8: Error Detail
    $ws.Range("P8").Value = $errorDetail
}

Update-HandbackRow "zh-cn" "a14ebcf6-3fca-4c76-925c-380c384446cc.1277613605aad44786cf4f78666a83ef22701133.zh-cn.xlf" "2016-08-30 18:54:14"
Update-HandbackRow "de-de" "a14ebcf6-3fca-4c76-925c-380c384446cc.1277613605aad44786cf4f78666a83ef22701133.de-de.xlf" "2016-08-30 18:54:22"
